$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 37039716
$ws.Range("I76").Value = 45457210
$ws.Range("J76").Value = 2740
$ws.Range("K76").Value = 45457210
$ws.Range("L76").Value = 2740
$ws.Range("M76").Value = -45456895
$ws.Range("N76").Value = -3370

$ws.Range("H79").Value = 37039716
$ws.Range("I79").Value = 45457210
$ws.Range("J79").Value = 2740
$ws.Range("K79").Value = 45457210
$ws.Range("L79").Value = 2740
$ws.Range("M79").Value = -45456118
$ws.Range("N79").Value = -4924

$ws.Range("H113").Value = 11113011
$ws.Range("I113").Value = 16668267
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 16668267
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = -16665013
$ws.Range("N113").Value = -9008

$ws.Range("H132").Value = 1823328.2
$ws.Range("I132").Value = 1563.14
$ws.Range("J132").Value = 10104079
$ws.Range("K132").Value = 4689.42
$ws.Range("L132").Value = 30312237
$ws.Range("M132").Value = -2159.42
$ws.Range("N132").Value = -30317297

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3466505
$ws.Range("I61").Value = 2084205.1
$ws.Range("J61").Value = 6538283
$ws.Range("K61").Value = 2084205.1
$ws.Range("L61").Value = 6538283
$ws.Range("M61").Value = -2083993.1
$ws.Range("N61").Value = -6538707

$ws.Range("H63").Value = 1942.0834
$ws.Range("I63").Value = 1700.7142
$ws.Range("J63").Value = 2280
$ws.Range("K63").Value = 1700.7142
$ws.Range("L63").Value = 2280
$ws.Range("M63").Value = -1014.7142
$ws.Range("N63").Value = -3652

$ws.Range("H66").Value = 1942.0834
$ws.Range("I66").Value = 1700.7142
$ws.Range("J66").Value = 2280
$ws.Range("K66").Value = 8503.571
$ws.Range("L66").Value = 11400
$ws.Range("M66").Value = -5071.571
$ws.Range("N66").Value = -18264

$ws.Range("H88").Value = 3672.64
$ws.Range("I88").Value = 2228
$ws.Range("J88").Value = 4234.4443
$ws.Range("K88").Value = 2228
$ws.Range("L88").Value = 4234.4443
$ws.Range("M88").Value = -1822
$ws.Range("N88").Value = -5046.4443

$ws.Range("H91").Value = 3672.64
$ws.Range("I91").Value = 2228
$ws.Range("J91").Value = 4234.4443
$ws.Range("K91").Value = 2228
$ws.Range("L91").Value = 4234.4443
$ws.Range("M91").Value = -824
$ws.Range("N91").Value = -7042.4443

$ws.Range("H97").Value = 636.3570999999999
$ws.Range("I97").Value = 640.9
$ws.Range("K97").Value = 640.9
$ws.Range("M97").Value = -144.9

$ws.Range("H132").Value = 23154120
$ws.Range("I132").Value = 29418836
$ws.Range("J132").Value = 7939807.5
$ws.Range("K132").Value = 88256508
$ws.Range("L132").Value = 23819422.5
$ws.Range("M132").Value = -88253978
$ws.Range("N132").Value = -23824482.5

$ws.Range("H136").Value = 3466505
$ws.Range("I136").Value = 2084205.1
$ws.Range("J136").Value = 6538283
$ws.Range("K136").Value = 6252615.300000001
$ws.Range("L136").Value = 19614849
$ws.Range("M136").Value = -6250065.300000001
$ws.Range("N136").Value = -19619949

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1902.34
$ws.Range("I86").Value = 1944.8989
$ws.Range("J86").Value = 1558
$ws.Range("K86").Value = 1944.8989
$ws.Range("L86").Value = 1558
$ws.Range("M86").Value = -821.8988999999999
$ws.Range("N86").Value = -3804

$ws.Range("H89").Value = 1902.34
$ws.Range("I89").Value = 1944.8989
$ws.Range("J89").Value = 1558
$ws.Range("K89").Value = 9724.494499999999
$ws.Range("L89").Value = 7790
$ws.Range("M89").Value = -4108.494499999999
$ws.Range("N89").Value = -19022

$ws.Range("H94").Value = 1878.4667
$ws.Range("I94").Value = 1252.0769
$ws.Range("J94").Value = 5950
$ws.Range("K94").Value = 1252.0769
$ws.Range("L94").Value = 5950
$ws.Range("M94").Value = -801.0769
$ws.Range("N94").Value = -6852

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1181935.5
$ws.Range("I31").Value = 965.9048
$ws.Range("K31").Value = 965.9048
$ws.Range("M31").Value = -670.9048

$ws.Range("H34").Value = 1181935.5
$ws.Range("I34").Value = 965.9048
$ws.Range("K34").Value = 965.9048
$ws.Range("M34").Value = -763.9048

$ws.Range("H62").Value = 22730140
$ws.Range("I62").Value = 2381.25
$ws.Range("J62").Value = 83337500
$ws.Range("K62").Value = 2381.25
$ws.Range("L62").Value = 83337500
$ws.Range("M62").Value = -1757.25
$ws.Range("N62").Value = -83338748

$ws.Range("H65").Value = 22730140
$ws.Range("I65").Value = 2381.25
$ws.Range("J65").Value = 83337500
$ws.Range("K65").Value = 11906.25
$ws.Range("L65").Value = 416687500
$ws.Range("M65").Value = -8786.25
$ws.Range("N65").Value = -416693740

$ws.Range("H134").Value = 741697.5
$ws.Range("I134").Value = 951.3
$ws.Range("K134").Value = 2853.9
$ws.Range("M134").Value = -318.8999999999996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 855.6818
$ws.Range("I86").Value = 727.3333
$ws.Range("J86").Value = 875.9474
$ws.Range("K86").Value = 2181.9999
$ws.Range("L86").Value = 2627.8422
$ws.Range("M86").Value = -995.9998999999998
$ws.Range("N86").Value = -4999.8422

$ws.Range("H89").Value = 855.6818
$ws.Range("I89").Value = 727.3333
$ws.Range("J89").Value = 875.9474
$ws.Range("K89").Value = 6545.9997
$ws.Range("L89").Value = 7883.5266
$ws.Range("M89").Value = -617.9997000000003
$ws.Range("N89").Value = -19739.5266

$ws.Range("H92").Value = 1470264.9
$ws.Range("I92").Value = 201.77777
$ws.Range("K92").Value = 605.33331
$ws.Range("M92").Value = 642.66669

$ws.Range("H107").Value = 1026468.8
$ws.Range("J107").Value = 1300.2667
$ws.Range("L107").Value = 3900.800099999999
$ws.Range("N107").Value = -7740.800099999999

$ws.Range("H122").Value = 1177.9333
$ws.Range("J122").Value = 1859.6471
$ws.Range("L122").Value = 16736.8239
$ws.Range("N122").Value = -21636.8239

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3707892.8
$ws.Range("I70").Value = 1920122.8
$ws.Range("J70").Value = 6948225.5
$ws.Range("K70").Value = 1920122.8
$ws.Range("L70").Value = 6948225.5
$ws.Range("M70").Value = -1919852.8
$ws.Range("N70").Value = -6948765.5

$ws.Range("H73").Value = 3707892.8
$ws.Range("I73").Value = 1920122.8
$ws.Range("J73").Value = 6948225.5
$ws.Range("K73").Value = 1920122.8
$ws.Range("L73").Value = 6948225.5
$ws.Range("M73").Value = -1919186.8
$ws.Range("N73").Value = -6950097.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1618.75
$ws.Range("I7").Value = 1991.6
$ws.Range("J7").Value = 997.3333
$ws.Range("K7").Value = 1991.6
$ws.Range("L7").Value = 997.3333
$ws.Range("M7").Value = -1879.6
$ws.Range("N7").Value = -1221.3333

$ws.Range("H68").Value = 2134.7144
$ws.Range("I68").Value = 1626.6666
$ws.Range("J68").Value = 2515.75
$ws.Range("K68").Value = 1626.6666
$ws.Range("L68").Value = 2515.75
$ws.Range("M68").Value = -877.6666
$ws.Range("N68").Value = -4013.75

$ws.Range("H71").Value = 2134.7144
$ws.Range("I71").Value = 1626.6666
$ws.Range("J71").Value = 2515.75
$ws.Range("K71").Value = 8133.333000000001
$ws.Range("L71").Value = 12578.75
$ws.Range("M71").Value = -4389.333000000001
$ws.Range("N71").Value = -20066.75

$ws.Range("H122").Value = 20118490
$ws.Range("I122").Value = 2130339.5
$ws.Range("K122").Value = 6391018.5
$ws.Range("M122").Value = -6388568.5

$ws.Range("H126").Value = 1618.75
$ws.Range("I126").Value = 1991.6
$ws.Range("J126").Value = 997.3333
$ws.Range("K126").Value = 5974.799999999999
$ws.Range("L126").Value = 2991.9999
$ws.Range("M126").Value = -3504.799999999999
$ws.Range("N126").Value = -7931.9999

$ws.Range("H132").Value = 2696513.5
$ws.Range("I132").Value = 3322955
$ws.Range("J132").Value = 2815.6
$ws.Range("K132").Value = 9968865
$ws.Range("L132").Value = 8446.799999999999
$ws.Range("M132").Value = -9966335
$ws.Range("N132").Value = -13506.8

$ws.Range("H136").Value = 2526261.8
$ws.Range("I136").Value = 2646512.2
$ws.Range("K136").Value = 7939536.600000001
$ws.Range("M136").Value = -7936986.600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1107.1807
$ws.Range("I136").Value = 388.10638
$ws.Range("J136").Value = 2045.9722
$ws.Range("K136").Value = 1164.31914
$ws.Range("L136").Value = 6137.9166
$ws.Range("M136").Value = 1385.68086
$ws.Range("N136").Value = -11237.9166
